$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taul1")

# Row 8: K8 gets value 2.1 (style already set to 3)
$ws.Cells.Item(8, 11).Value = 2.1

# Row 12: K12 value removed (cell cleared back to empty, but keep style s="3")
$ws.Cells.Item(12, 11).ClearContents()

# Row 22: K22 value changes from 2.65 to 1.8
$ws.Cells.Item(22, 11).Value = 1.8

# Row 25: K25 value removed; L25 value changes
$ws.Cells.Item(25, 11).ClearContents()
$ws.Cells.Item(25, 12).Value = 1.1620900000000001

# Row 26: K26 value changes
$ws.Cells.Item(26, 11).Value = 0.83333000000000002

# Row 27: K27 gets a value; L27 value changes
$ws.Cells.Item(27, 11).Value = 1.1666700000000001
$ws.Cells.Item(27, 12).Value = 1.0040100000000001

# Row 32: A32 "Jerry Miculek" -> "Ben Stoeger"; B32 value changes; G32 gets empty styled cell
$ws.Cells.Item(32, 1).Value = "Ben Stoeger"
$ws.Cells.Item(32, 2).Value = 0.86399999999999999
$ws.Cells.Item(32, 7).NumberFormat = "0.00%"

# Row 33: A33 "Ben Stoeger" -> "Jerry Miculek"; B33 value changes
$ws.Cells.Item(33, 1).Value = "Jerry Miculek"
$ws.Cells.Item(33, 2).Value = 0.75

# Sheet view: active cell / selection change
$ws.Range("F32:G32").Select()
